# Apply "week4 kaggle baseline beaten" edit:
# - Move the "model" column (header + LogisticRegression value) from column J to column K
# - Insert a new "scaler" header in J6
# - Add a new row 8 with GradientBoostingClassifier results (incl. StandardScaler)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "model" column from J to K
$ws.Range("K5").Value = $ws.Range("J5").Value2
$ws.Range("J5").Value = $null

$ws.Range("K6").Value = $ws.Range("J6").Value2
$ws.Range("J6").Value = "scaler"

$ws.Range("K7").Value = $ws.Range("J7").Value2
$ws.Range("J7").Value = $null

# New row 8: GradientBoostingClassifier results.
# The numeric-looking metrics must land as shared-string TEXT (matching the
# source file), not as numbers, so force text format before entry then
# strip the formatting override back off (keeps them type "s" with the
# default/implicit cell style, same as the rest of the sheet).
$numericTextRange = $ws.Range("B8:G8")
$numericTextRange.NumberFormat = "@"
$ws.Range("B8").Value = "0.7090677210791867"
$ws.Range("C8").Value = "0.47916666666666663"
$ws.Range("D8").Value = "0.0024509803921568627"
$ws.Range("E8").Value = "0.004871442393976412"
$ws.Range("F8").Value = "0.7127077471277219"
$ws.Range("G8").Value = "0.67713"
$numericTextRange.ClearFormats()

$ws.Range("H8").Value = "most_frequent"
$ws.Range("I8").Value = "OrdinalEncoder"
$ws.Range("J8").Value = "StandardScaler"
$ws.Range("K8").Value = "GradientBoostingClassifier"

$ws.Range("L8").Select()
